$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The value that used to live in C14 (2 hours, shown with the [hh]:mm:ss
# style already applied to C13/C14) moves up into C13, and C14 is cleared
# out. The stray "Time Spent " label that used to sit in C13 goes away
# entirely (it was never referenced from anywhere else, so Excel will drop
# it from the shared-strings table on save once nothing uses it).
$ws.Range("C13").Value = 0.083333333333333329
$ws.Range("C14").ClearContents()

# Move the active selection onto C14, matching the saved selection state.
$ws.Range("C14").Select()
